$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap test case content between row 2 block (Pos_Fun_0001) and row 14 block (Pos_Fun_0004) ---
# Row 2 becomes the former row-14 "Convert question" test case content
$ws.Range("B2").Value = "Convert question"
$ws.Range("D2").Value = "oyaata kohomadha?"
$ws.Range("E2").Value = "ඔයාට කොහොමද?"
$ws.Range("F2").Value = "ඔයාට කොහොමද?"
$ws.Range("H2").Value = "Interrogative tone identified."
$ws.Range("I2").Value = "Interrogative (question)"
$ws.Range("H3").Value = "Question mark retained at end."
$ws.Range("I3").Value = "Question mark usage"
$ws.Range("H4").Value = "Correct word choice for `"kohomadha`"."

# Row 14 becomes the former row-2 "Convert simple daily sentence" test case content
$ws.Range("B14").Value = "Convert simple daily sentence"
$ws.Range("D14").Value = "api heta beach yamu"
$ws.Range("E14").Value = "අපි හෙට beach යමු"
$ws.Range("F14").Value = "අපි හෙට beach යමු"
$ws.Range("H14").Value = "The sentence meaning is correctly preserved."
$ws.Range("I14").Value = "Daily language usage"
$ws.Range("H15").Value = "Sinhala spelling is accurate."
$ws.Range("I15").Value = "Simple sentence"
$ws.Range("H16").Value = "Sentence structure is correctly converted."

# --- Update sheet view: zoom to 70% and change the active selection to J6 ---
$ws.Activate()
$ws.Range("J6").Select() | Out-Null
$excel.ActiveWindow.Zoom = 70

